$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function ReplaceInCell($row, $col, $old, $new) {
    $c = $t.Cell($row, $col)
    $rng = $d.Range($c.Range.Start, $c.Range.End)
    $res = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1)
    return $res
}

# Row 1: "100" -> "0M"
ReplaceInCell 1 1 "100" "0M" | Out-Null

# Row 2: "0" -> "0M"
ReplaceInCell 2 1 "0" "0M" | Out-Null

# Row 3: "53" -> "0M"
ReplaceInCell 3 1 "53" "0M" | Out-Null

# Row 4: "3" -> "43"
ReplaceInCell 4 1 "3" "43" | Out-Null

# Row 6: "0.00004" -> "0.00005"
ReplaceInCell 6 1 "0.00004" "0.00005" | Out-Null

# Row 10: "0.00003" -> "0.00004"
ReplaceInCell 10 1 "0.00003" "0.00004" | Out-Null

# Row 12: "0.00010" -> "0.00149"
ReplaceInCell 12 1 "0.00010" "0.00149" | Out-Null

# Row 44: collapse the multi-run "5 <tab> 0.00003 <tab> ... 100.0" into just "100"
$t.Cell(44,1).Range.Text = "100"

# Row 45: collapse the multi-run "25 <tab> 0.00003 <tab> ... 100.0" into just "0"
$t.Cell(45,1).Range.Text = "0"

# Row 46: collapse the multi-run "10 <tab> 0.00003 <tab> ... 100.0" into just "53"
$t.Cell(46,1).Range.Text = "53"
